$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 19704
$ws1.Range("F5").Value = 792
$ws1.Range("F8").Value = 11
$ws1.Range("F9").Value = 7457
$ws1.Range("F10").Value = 496
$ws1.Range("F14").Value = 149
$ws1.Range("F16").Value = 2
$ws1.Range("F20").Value = 386
$ws1.Range("F23").Value = 46
$ws1.Range("F27").Value = 1080
$ws1.Range("F29").Value = 14
$ws1.Range("F30").Value = 171
$ws1.Range("F34").Value = 2787
$ws1.Range("F37").Value = 17
$ws1.Range("F38").Value = 12548
$ws1.Range("F44").Value = 348
$ws1.Range("F47").Value = 93

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 33

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 19704
$ws4.Range("F5").Value = 0
$ws4.Range("F9").Value = 7457
$ws4.Range("F10").Value = 496
$ws4.Range("F15").Value = 106
$ws4.Range("F17").Value = 231
$ws4.Range("F19").Value = 1331
$ws4.Range("F30").Value = 171
$ws4.Range("F31").Value = 5220
$ws4.Range("F32").Value = 556
$ws4.Range("F34").Value = 50
$ws4.Range("F36").Value = 2787
$ws4.Range("F39").Value = 17
$ws4.Range("F40").Value = 12548
$ws4.Range("F49").Value = 93
